$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Candidate ID 231004201 -> 231011130, plus new credential strings)
$ws.Range("A2").Value = 'gudbN200'
$ws.Range("B2").Value = 231011130
$ws.Range("C2").Value = 'cptaakv72'
$ws.Range("D2").Value = 'x#$CJ6j2'
$ws.Range("F2").Value = 'yseVZrQU'
$ws.Range("G2").Value = 'xeGa'

# Row 3 (Candidate ID 231004200 -> 231011129, plus new credential strings)
$ws.Range("A3").Value = 'mzlLH609'
$ws.Range("B3").Value = 231011129
$ws.Range("C3").Value = 'uuwozii72'
$ws.Range("D3").Value = 'u5U%&aK3'
$ws.Range("F3").Value = 'sfcHflqi'
$ws.Range("G3").Value = 'iSIt'
